$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 45210
$ws.Range("J2").Value = 550
$ws.Range("K2").Value = 1500
$ws.Range("L2").Value = 1600
$ws.Range("M2").Value = 1536
$ws.Range("P2").Value = 1536

$ws.Range("D3").Value = 44907
$ws.Range("J3").Value = 2300
$ws.Range("K3").Value = 900
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = 952
$ws.Range("P3").Value = 952

$ws.Range("D4").Value = 45132
$ws.Range("J4").Value = 170
$ws.Range("K4").Value = 2200
$ws.Range("L4").Value = 2500
$ws.Range("M4").Value = 2359
$ws.Range("P4").Value = 2359

$ws.Range("D5").Value = 44895
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 1200
$ws.Range("L5").Value = 1300
$ws.Range("M5").Value = 1255
$ws.Range("P5").Value = 1255

$ws.Range("D6").Value = 44537
$ws.Range("J6").Value = 800
$ws.Range("K6").Value = 1300
$ws.Range("L6").Value = 1400
$ws.Range("M6").Value = 1350
$ws.Range("P6").Value = 1350

$ws.Range("D7").Value = 45204
$ws.Range("J7").Value = 1200
$ws.Range("K7").Value = 1600
$ws.Range("L7").Value = 1700
$ws.Range("M7").Value = 1650
$ws.Range("P7").Value = 1650

$ws.Range("D8").Value = 45212
$ws.Range("J8").Value = 750
$ws.Range("K8").Value = 1400
$ws.Range("L8").Value = 1500
$ws.Range("M8").Value = 1440
$ws.Range("P8").Value = 1440

$ws.Range("D9").Value = 45062
$ws.Range("J9").Value = 1700
$ws.Range("K9").Value = 2800
$ws.Range("L9").Value = 3000
$ws.Range("M9").Value = 2900
$ws.Range("P9").Value = 2900

$ws.Range("D10").Value = 44200
$ws.Range("J10").Value = 1500
$ws.Range("K10").Value = 1400
$ws.Range("L10").Value = 1500
$ws.Range("M10").Value = 1450
$ws.Range("P10").Value = 1450

$ws.Range("D11").Value = 44210
$ws.Range("J11").Value = 1450
$ws.Range("K11").Value = 1600
$ws.Range("L11").Value = 1700
$ws.Range("M11").Value = 1650
$ws.Range("P11").Value = 1650

$ws.Range("D12").Value = 45205
$ws.Range("J12").Value = 3500
$ws.Range("K12").Value = 1400
$ws.Range("L12").Value = 1500
$ws.Range("M12").Value = 1457
$ws.Range("P12").Value = 1457

$ws.Range("D13").Value = 44638
$ws.Range("J13").Value = 800
$ws.Range("K13").Value = 2500
$ws.Range("L13").Value = 2800
$ws.Range("M13").Value = 2650
$ws.Range("P13").Value = 2650

$ws.Range("D14").Value = 44175
$ws.Range("J14").Value = 1400
$ws.Range("K14").Value = 1900
$ws.Range("L14").Value = 2000
$ws.Range("M14").Value = 1950
$ws.Range("P14").Value = 1950

$ws.Range("D15").Value = 45203
$ws.Range("J15").Value = 800
$ws.Range("K15").Value = 1800
$ws.Range("L15").Value = 2000
$ws.Range("M15").Value = 1900
$ws.Range("P15").Value = 1900

$ws.Range("D16").Value = 44883
$ws.Range("J16").Value = 290
$ws.Range("K16").Value = 1400
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = 1434
$ws.Range("P16").Value = 1434

$ws.Range("D17").Value = 44893
$ws.Range("J17").Value = 3300
$ws.Range("K17").Value = 1200
$ws.Range("L17").Value = 1300
$ws.Range("M17").Value = 1261
$ws.Range("P17").Value = 1261
